# Update Emily Meyer's information in the People spreadsheet.
# - Remove the old, minimal "Meyer, Emily" row (row 3), leaving it blank
#   like the other spacer rows.
# - Re-add Emily Meyer as a full entry at the bottom of the table (row 30)
#   with her name/email plus her new Title, Description and LinkedIn link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember Emily's name/email (currently sitting in row 3) before we touch it.
$emilyName  = $ws.Range("B3").Value()
$emilyEmail = $ws.Range("C3").Value()

# --- Clear the old row 3 entry so it becomes an empty spacer row -----------
$ws.Rows.Item(3).Clear()
$ws.Rows.Item(3).AutoFit()

# --- Build the new row 30 ---------------------------------------------------
# Copy the formatting of an existing, fully-populated data row (row 7: Name,
# Email, Title, Link, Description) onto row 30 so fonts/fills/number formats
# match the rest of the sheet.
$ws.Range("B7:K7").Copy()
$ws.Range("B30:K30").PasteSpecial(-4122)

# Column E's highlight format on row 7 is slightly different from the other
# highlighted cells in that row (D/G/H/I/K) -- match those instead, as the
# target layout does.
$ws.Range("D29").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(30).RowHeight = $ws.Rows.Item(7).RowHeight

# Fill in the values.
$ws.Range("B30").Value = $emilyName
$ws.Range("C30").Value = $emilyEmail
$ws.Range("E30").Value = "Post Bac Researcher"
$ws.Range("J30").Value = "I am a postbac fellow at the National Institutes of Health. I am interested in visual neuroscience, specifically comparing behavior with higher-level processing through experimental and computational methods. I graduated Tulane in 2019 with majors in mathematics and neuroscience and a public health minor. With the Mathematical Modeling and Analysis Lab, I studied the effects of diabetes on the progression and drug resistance of tuberculosis. I also previously worked on modeling vector-borne diseases with multiple risk groups and behavioral factors."

# Add her LinkedIn link as a real hyperlink.
$ws.Hyperlinks.Add($ws.Range("F30"), "https://www.linkedin.com/in/emily-meyer-769a95158/") | Out-Null

# Adding the hyperlink re-applies Excel's default Hyperlink font; restore the
# cell's format to match the rest of the highlighted link cells (style used
# by F7, the other Link-1 cell).
$ws.Range("F7").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view: scroll down and select the newly added row ----------
$ws.Range("B30:K30").Select()
